# Auto commit on 02-2025-time
# Extends the "Report" sheet from 8 data rows to 12 rows (4 new work orders),
# widens the print area to match, and nudges the active selection the same
# way Excel leaves it after someone fills data further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# --- 1. Clone row formatting for the 4 new rows before writing values -----
# Row 7 is the "odd" banding (style family 7/8/9), row 8 is the "even"
# banding (style family 3/4/4). New rows 9 and 11 follow the odd banding,
# rows 10 and 12 follow the even banding.
$ws.Range("A7:AK7").Copy()
$ws.Range("A9:AK9").PasteSpecial(-4122)
$ws.Range("A7:AK7").Copy()
$ws.Range("A11:AK11").PasteSpecial(-4122)

$ws.Range("A8:AK8").Copy()
$ws.Range("A10:AK10").PasteSpecial(-4122)
$ws.Range("A8:AK8").Copy()
$ws.Range("A12:AK12").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# Rows 9/11 (odd banding) and row 10 (even banding) all pick up the
# wrap-on variant of the P/AC (work-content) style, matching P7/AC7. The
# base format paste above already carried the rest of row 7/8's look, so
# just flip WrapText on for those four cells (PasteSpecial doesn't carry
# WrapText through, so it has to be set explicitly).
$ws.Range("P9").WrapText = $true
$ws.Range("AC9").WrapText = $true
$ws.Range("P10").WrapText = $true
$ws.Range("AC10").WrapText = $true
$ws.Range("P11").WrapText = $true
$ws.Range("AC11").WrapText = $true

# --- 2. Row 9 values --------------------------------------------------------
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "服務"
$ws.Range("C9").Value = 2025070395
$ws.Range("F9").Value = 4191
$ws.Range("G9").Value = "三重溪美店"
$ws.Range("H9").Value = "新北市三重區"
$ws.Range("Q9").Value = "THILF04191"
$ws.Range("R9").Value = "新北一"
$ws.Range("S9").Value = "吳宗鴻"
$ws.Range("T9").Value = 1
$ws.Range("U9").Value = "已完工"
$ws.Range("V9").Value = "2025-07-02 15:04:28"
$ws.Range("W9").Value = "2025-07-02 14:00:00"
$ws.Range("X9").Value = "2025-07-02 15:03:00"
$ws.Range("Z9").Value = 1.1
$ws.Range("AB9").Value = "到場處理"
$ws.Range("AC9").Value = "PMQ3+TVV+STAR"
$ws.Range("AD9").Value = "O"
$ws.Range("AJ9").Value = "O"
$ws.Range("AK9").Value = "O"

# --- 3. Row 10 values -------------------------------------------------------
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "服務"
$ws.Range("C10").Value = 2025070417
$ws.Range("F10").Value = 4210
$ws.Range("G10").Value = "三重福華店"
$ws.Range("H10").Value = "新北市三重區"
$ws.Range("Q10").Value = "THILF04210"
$ws.Range("R10").Value = "新北一"
$ws.Range("S10").Value = "吳宗鴻"
$ws.Range("T10").Value = 1
$ws.Range("U10").Value = "已完工"
$ws.Range("V10").Value = "2025-07-02 16:09:45"
$ws.Range("W10").Value = "2025-07-02 15:10:00"
$ws.Range("X10").Value = "2025-07-02 16:08:00"
$ws.Range("Z10").Value = 1
$ws.Range("AB10").Value = "到場處理"
$ws.Range("AC10").Value = "PMQ3+STAR"
$ws.Range("AD10").Value = "O"
$ws.Range("AJ10").Value = "O"
$ws.Range("AK10").Value = "O"

# --- 4. Row 11 values -------------------------------------------------------
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "服務"
$ws.Range("C11").Value = 2025070428
$ws.Range("F11").Value = 4210
$ws.Range("G11").Value = "三重福華店"
$ws.Range("H11").Value = "新北市三重區"
$ws.Range("Q11").Value = "THILF04210"
$ws.Range("R11").Value = "新北一"
$ws.Range("S11").Value = "吳宗鴻"
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = "已完工"
$ws.Range("V11").Value = "2025-07-02 17:13:17"
$ws.Range("W11").Value = "2025-07-02 15:10:00"
$ws.Range("X11").Value = "2025-07-02 16:08:00"
$ws.Range("Z11").Value = 1
$ws.Range("AB11").Value = "到場處理"
$ws.Range("AC11").Value = "PMQ3+STAR"
$ws.Range("AD11").Value = "O"
$ws.Range("AJ11").Value = "O"
$ws.Range("AK11").Value = "O"

# --- 5. Row 12 values -------------------------------------------------------
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "服務"
$ws.Range("C12").Value = 2025070429
$ws.Range("F12").Value = "D194"
$ws.Range("G12").Value = "北縣五華三店"
$ws.Range("H12").Value = "新北市三重區"
$ws.Range("Q12").Value = "THILF0D194"
$ws.Range("R12").Value = "新北一"
$ws.Range("S12").Value = "吳宗鴻"
$ws.Range("T12").Value = 1
$ws.Range("U12").Value = "已完工"
$ws.Range("V12").Value = "2025-07-02 17:13:54"
$ws.Range("W12").Value = "2025-07-02 16:20:00"
$ws.Range("X12").Value = "2025-07-02 17:13:00"
$ws.Range("Z12").Value = 0.9
$ws.Range("AB12").Value = "到場處理"
$ws.Range("AC12").Value = "PMQ3+STAR"
$ws.Range("AD12").Value = "O"
$ws.Range("AJ12").Value = "O"
$ws.Range("AK12").Value = "O"

# --- 6. Row 8's own P/AC cells switch to the wrap-on style too -------------
$ws.Range("P8").WrapText = $true
$ws.Range("AC8").WrapText = $true

# --- 7. Print area + active selection --------------------------------------
$ws.PageSetup.PrintArea = '$A$1:$AK$12'
$ws.Range("AC9").Select()
